$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 249
$wsOff.Range("C3").Value = 156
$wsOff.Range("D3").Value = 50
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 6

# DEF sheet - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 261
$wsDef.Range("C3").Value = 196
$wsDef.Range("D3").Value = 60
$wsDef.Range("E3").Value = 33
